$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 currently holds the text "R40". The rule's "From" label changes
# to the text "1" (still a text value, not a number) -- force text so the
# numeric-looking string isn't auto-converted to a number.
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
